{"js": "// Office.js (Word JavaScript API) edit script\n// Body of: async (context) => { ... }\n\n// ---------------------------------------------------------------------\n// 1) Insert a new \"Meta description\" paragraph right after the H1 title\n//    (\"Play Black Widow Slot for Free - Pros, Cons & Comparison\").\n// ---------------------------------------------------------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\n\nconst metaPara = titlePara.insertParagraph(\"\", Word.InsertLocation.after);\n// Match the body-text look (no explicit heading style) used elsewhere in\n// the document; assigning \"Normal\" keeps the emitted markup free of an\n// explicit <w:pPr><w:pStyle/> (Normal is the implicit default style).\nmetaPara.style = \"Normal\";\n\nconst labelRun = metaPara.insertText(\"Meta description\", Word.InsertLocation.end);\nlabelRun.font.bold = true;\n\nconst restRun = metaPara.insertText(\n  \": Read our Black Widow slot review and play for free! Discover pros and cons, bonus features, and a comparison with similar slots.\",\n  Word.InsertLocation.end\n);\n// Explicitly clear bold on the second run so it doesn't inherit the bold\n// formatting of the insertion point (keeps it as a separate, un-bolded run).\nrestRun.font.bold = false;\n\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Remove the duplicated bold title paragraph near the end of the\n//    document, and rewrite the remaining italic paragraph's text into an\n//    image-generation prompt (keeping its italic formatting).\n// ---------------------------------------------------------------------\nconst allParagraphs = context.document.body.paragraphs;\nallParagraphs.load(\"text\");\nawait context.sync();\n\nconst count = allParagraphs.items.length;\nconst duplicateTitlePara = allParagraphs.items[count - 2];\nduplicateTitlePara.delete();\nawait context.sync();\n\nconst refreshedParagraphs = context.document.body.paragraphs;\nrefreshedParagraphs.load(\"items\");\nawait context.sync();\n\nconst newCount = refreshedParagraphs.items.length;\nconst descriptionPara = refreshedParagraphs.items[newCount - 1];\n\nconst oldTextResults = descriptionPara.search(\n  \"Read our Black Widow slot review and play for free! Discover pros and cons, bonus features, and a comparison with similar slots.\",\n  { matchCase: true }\n);\nawait context.sync();\n\noldTextResults.items[0].insertText(\n  \"Create a feature image for Black Widow slot game featuring a happy Maya warrior with glasses in a cartoon style. The image should feature the warrior holding a spider and standing in front of a spider web. The background should be dark with cobwebs on the corners to match the theme of the game. The warrior should be dressed in a black jumpsuit with a red hourglass symbol on the chest and his/her arms folded in front. The image should be eye-catching with vibrant colors to attract players to the game.\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Insert a new \"Meta description\" paragraph right after the H1 title\n#    (\"Play Black Widow Slot for Free - Pros, Cons & Comparison\").\n# ---------------------------------------------------------------------\n$titlePara = $d.Paragraphs(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs(2)\n# Match the body-text look used elsewhere in the document (no explicit\n# heading style); \"Normal\" is the implicit default so no <w:pPr><w:pStyle>\n# gets written out.\n$metaPara.Style = \"Normal\"\n\n$metaRange = $metaPara.Range\n$labelText = \"Meta description\"\n$metaRange.Text = $labelText\n\n# Bold only the literal text (not the paragraph mark) so the paragraph's\n# mark-run-properties stay untouched.\n$boldRange = $d.Range($metaRange.Start, $metaRange.Start + $labelText.Length)\n$boldRange.Font.Bold = 1\n\n# Append the remainder of the sentence as its own (non-bold) run.\n$metaPara2 = $d.Paragraphs(2)\n$insertPoint = $d.Range($metaPara2.Range.End - 1, $metaPara2.Range.End - 1)\n$insertPoint.InsertAfter(\": Read our Black Widow slot review and play for free! Discover pros and cons, bonus features, and a comparison with similar slots.\")\n$insertPoint.Font.Bold = 0\n\n# ---------------------------------------------------------------------\n# 2) Remove the duplicated bold title paragraph near the end of the\n#    document, and rewrite the remaining italic paragraph's text into an\n#    image-generation prompt (keeping its italic formatting).\n# ---------------------------------------------------------------------\n$count = $d.Paragraphs.Count\n$duplicateTitlePara = $d.Paragraphs($count - 1)\n$duplicateTitlePara.Range.Delete()\n\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$findRange = $lastPara.Range\n$findRange.Find.Execute(\n    \"Read our Black Widow slot review and play for free! Discover pros and cons, bonus features, and a comparison with similar slots.\",\n    $false, $true, $false, $false, $false, $true, 1, $false,\n    \"Create a feature image for Black Widow slot game featuring a happy Maya warrior with glasses in a cartoon style. The image should feature the warrior holding a spider and standing in front of a spider web. The background should be dark with cobwebs on the corners to match the theme of the game. The warrior should be dressed in a black jumpsuit with a red hourglass symbol on the chest and his/her arms folded in front. The image should be eye-catching with vibrant colors to attract players to the game.\",\n    2\n)\n"}
